$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2292682926829268
$ws.Range("C2").Value = 0.4878048780487805
$ws.Range("J2").Value = 0.00975609756097561
$ws.Range("P2").Value = 0.1707317073170732
$ws.Range("S2").Value = 0.1024390243902439
$ws.Range("B3").Value = 0.03846153846153846
$ws.Range("C3").Value = 0.02884615384615385
$ws.Range("J3").Value = 0.03846153846153846
$ws.Range("P3").Value = 0.7211538461538461
$ws.Range("S3").Value = 0.1730769230769231
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("P4").Value = 0.6136363636363636
$ws.Range("S4").Value = 0.3181818181818182
$ws.Range("B6").Value = 0.08426966292134831
$ws.Range("D6").Value = 0.01123595505617977
$ws.Range("F6").Value = 0.05617977528089887
$ws.Range("J6").Value = 0.1685393258426966
$ws.Range("Q6").Value = 0.1573033707865168
$ws.Range("R6").Value = 0.07865168539325842
$ws.Range("S6").Value = 0.4438202247191011
$ws.Range("B7").Value = 0.04458598726114649
$ws.Range("D7").Value = 0.03184713375796178
$ws.Range("F7").Value = 0.05732484076433121
$ws.Range("J7").Value = 0.1592356687898089
$ws.Range("O7").Value = 0.006369426751592357
$ws.Range("Q7").Value = 0.1464968152866242
$ws.Range("R7").Value = 0.1146496815286624
$ws.Range("S7").Value = 0.4394904458598726
$ws.Range("B8").Value = 0.05833333333333333
$ws.Range("D8").Value = 0.01388888888888889
$ws.Range("F8").Value = 0.05833333333333333
$ws.Range("J8").Value = 0.1611111111111111
$ws.Range("O8").Value = 0.025
$ws.Range("Q8").Value = 0.1694444444444445
$ws.Range("R8").Value = 0.09444444444444444
$ws.Range("S8").Value = 0.4194444444444445
$ws.Range("B9").Value = 0.09876543209876543
$ws.Range("D9").Value = 0.01851851851851852
$ws.Range("F9").Value = 0.0308641975308642
$ws.Range("J9").Value = 0.1481481481481481
$ws.Range("O9").Value = 0.02469135802469136
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.09259259259259259
$ws.Range("S9").Value = 0.4197530864197531
$ws.Range("B10").Value = 0.08615948670944087
$ws.Range("D10").Value = 0.02841429880843263
$ws.Range("E10").Value = 0.0009165902841429881
$ws.Range("F10").Value = 0.08799266727772685
$ws.Range("J10").Value = 0.1237396883593034
$ws.Range("O10").Value = 0.005499541704857928
$ws.Range("Q10").Value = 0.1943171402383135
$ws.Range("R10").Value = 0.0999083409715857
$ws.Range("S10").Value = 0.3730522456461962
$ws.Range("G11").Value = 0.1477272727272727
$ws.Range("J11").Value = 0.08712121212121213
$ws.Range("K11").Value = 0.2083333333333333
$ws.Range("L11").Value = 0.553030303030303
$ws.Range("S11").Value = 0.003787878787878788
$ws.Range("G12").Value = 0.7210884353741497
$ws.Range("J12").Value = 0.2312925170068027
$ws.Range("L12").Value = 0.0272108843537415
$ws.Range("S12").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.5263157894736842
$ws.Range("J13").Value = 0.3421052631578947
$ws.Range("S13").Value = 0.131578947368421
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.6666666666666666
$ws.Range("F15").Value = 0.02824858757062147
$ws.Range("H15").Value = 0.1412429378531073
$ws.Range("I15").Value = 0.096045197740113
$ws.Range("J15").Value = 0.4406779661016949
$ws.Range("K15").Value = 0.05084745762711865
$ws.Range("M15").Value = 0.005649717514124294
$ws.Range("O15").Value = 0.06214689265536723
$ws.Range("S15").Value = 0.1751412429378531
$ws.Range("F16").Value = 0.007518796992481203
$ws.Range("H16").Value = 0.1879699248120301
$ws.Range("I16").Value = 0.06015037593984962
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("K16").Value = 0.1578947368421053
$ws.Range("M16").Value = 0.01503759398496241
$ws.Range("O16").Value = 0.06766917293233082
$ws.Range("S16").Value = 0.07518796992481203
$ws.Range("F17").Value = 0.0056657223796034
$ws.Range("H17").Value = 0.1614730878186969
$ws.Range("I17").Value = 0.09348441926345609
$ws.Range("J17").Value = 0.4645892351274787
$ws.Range("K17").Value = 0.1161473087818697
$ws.Range("M17").Value = 0.0141643059490085
$ws.Range("O17").Value = 0.05382436260623229
$ws.Range("S17").Value = 0.0906515580736544
$ws.Range("F18").Value = 0.01041666666666667
$ws.Range("H18").Value = 0.2083333333333333
$ws.Range("I18").Value = 0.078125
$ws.Range("J18").Value = 0.3802083333333333
$ws.Range("K18").Value = 0.08333333333333333
$ws.Range("M18").Value = 0.02083333333333333
$ws.Range("O18").Value = 0.09895833333333333
$ws.Range("S18").Value = 0.1197916666666667
$ws.Range("F19").Value = 0.009689922480620155
$ws.Range("H19").Value = 0.2054263565891473
$ws.Range("I19").Value = 0.08624031007751938
$ws.Range("J19").Value = 0.3701550387596899
$ws.Range("K19").Value = 0.1143410852713178
$ws.Range("M19").Value = 0.02616279069767442
$ws.Range("N19").Value = 0.002906976744186046
$ws.Range("O19").Value = 0.07267441860465117
$ws.Range("S19").Value = 0.1124031007751938
